$d = $word.ActiveDocument

# The last two reference paragraphs ("van Mourik..." and "Leslie, H. A...")
# are being removed, keeping only the "Capinha, L..." paragraph before
# the section break.
$startPara = $d.Paragraphs(4)
$endPara   = $d.Paragraphs(5)

$deleteRange = $d.Range($startPara.Range.Start, $endPara.Range.End)
$deleteRange.Delete()
